# CEDS Data and exogenous assumptions updates
# Insert a new sector row ("1A1bc_Other-feedstocks") into the "Sectors" sheet
# immediately before the existing "1A2a_Ind-Comb-Iron-steel" row (row 6),
# pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# Insert a new blank row at row 6 (shifts rows 6:60 down to 7:61)
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row with the new sector's data
$ws.Cells.Item(6, 1).Value = "1A1bc_Other-feedstocks"
$ws.Cells.Item(6, 2).Value = "Energy_Combustion"
$ws.Cells.Item(6, 3).Value = "kt"
$ws.Cells.Item(6, 4).Value = "NC"

# Match the saved file's reported selection/active cell (A6 on Sectors tab)
$ws.Range("A6").Select()
